$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New record row (row 2) - identifier, title, levelOfDescription, extentAndMedium, notes
# alternativeIdentifiers (B) / date_s (D) / file_path (H) are not populated for this
# record; D2 and H2 still pick up the row's formatting even though blank.
$cells = @("A2", "C2", "D2", "E2", "F2", "G2", "H2")
foreach ($addr in $cells) {
    $r = $ws.Range($addr)
    $r.Font.Name = "Calibri"
    $r.Font.ThemeColor = 1
}

$ws.Range("A2").Value = "MCH192"
$ws.Range("C2").Value = "ALBUM RE PACT OF SOLIDARITY BETWEEN THE ANC AND THE CITY OF REGGION EMIBA"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

$ws.Range("A2").Select()
$ws.Range("A2:H2").Select()
